$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Rename the "_old" / "_new" header-suffixes to "_FV2310" / "_FV2404"
#    (row 1, columns A:U). Using Find/Replace keeps this limited to the
#    header cells since no other cell in the sheet contains these
#    substrings.
# -----------------------------------------------------------------
$ws.Cells.Replace("_old", "_FV2310")
$ws.Cells.Replace("_new", "_FV2404")

# -----------------------------------------------------------------
# 2) Freeze the header row (split/freeze after row 1).
# -----------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# -----------------------------------------------------------------
# 3) Turn A1:U59 into an Excel Table ("Table1") without Excel
#    capturing the pre-existing bold/filled header formatting as a
#    table "headerRowDxf" override (which the source workbook does
#    not have). We do this by temporarily stashing the header's
#    formatting, clearing direct formatting from the header row
#    before the table is created, then pasting the formatting back
#    (as a single formats-only paste, so no extra/unused style
#    records get left behind in styles.xml).
# -----------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$stashCell = $ws.Range("AA1")

$ws.Range("A1").Copy()
$stashCell.PasteSpecial(-4122) # xlPasteFormats

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

$stashCell.Copy()
$headerRange.PasteSpecial(-4122) # xlPasteFormats

$stashCell.Clear()
$ws.Range("A1").Select()
